# Add a new "PassengerDetails" worksheet after the existing "SearchFlight" sheet
# and populate it with passenger utility data (firstname/lastname/mobilenumber/email),
# including a mailto: hyperlink on the email cell.

$wb = $excel.ActiveWorkbook

$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$newSheet = $wb.Worksheets.Add([System.Type]::Missing, $lastSheet)
$newSheet.Name = "PassengerDetails"

# Header row
$newSheet.Range("A1").Value = "firstname"
$newSheet.Range("C1").Value = "lastname"
$newSheet.Range("E1").Value = "mobilenumber"
$newSheet.Range("G1").Value = "email"

# Data row
$newSheet.Range("A2").Value = "Akash"
$newSheet.Range("C2").Value = "Kumar"
$newSheet.Range("E2").Value = 9876784563
$newSheet.Range("G2").Value = "akash@gmail.com"

# Hyperlink the email address
$newSheet.Hyperlinks.Add($newSheet.Range("G2"), "mailto:akash@gmail.com")

# Column E (mobile number) needs to widen to fit the long number, like Excel
# does automatically when a number overflows the default column width.
$newSheet.Columns.Item(5).EntireColumn.AutoFit()

# Make the new sheet the active / selected tab
$newSheet.Select() | Out-Null
$newSheet.Range("E3").Select() | Out-Null
